$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row correct-answer count (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update the "Total" row correct-answer count (B12: 33 -> 55)
$ws.Range("B12").Value = 55

# Update the correct/total marks summary text (E12: "22/84" -> "55/140")
$ws.Range("E12").Value = "55/140"
